# Apply the COC04 price-list update:
#  - refresh several product weight descriptions (e.g. "53G" -> "50G")
#  - remove the discontinued "KIT KAT CHO 2F 17G" (code 10036982) line,
#    shifting all subsequent rows up by one

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update product descriptions (column B) with new pack weights ---
$ws.Range("B6").Value  = "S/Q CHOCO CRISPY 50G"
$ws.Range("B7").Value  = "S/Q CHO FRUIT&NUT 52"
$ws.Range("B10").Value = "S/Q CHUNKY WHT 82G"
$ws.Range("B13").Value = "S/Q CHUNKY ALMOND 82"
$ws.Range("B14").Value = "S/Q CHOCO ALMOND 52G"
$ws.Range("B15").Value = "S/Q CHOCO CASHEW 52G"

$ws.Range("B45").Value = "FISH.EXTRA STRG PT22"
$ws.Range("B46").Value = "FISH.STRNG.MINT HJ22"
$ws.Range("B47").Value = "FISHERMAN'S LEMON 22"
$ws.Range("B48").Value = "FSHERMAN'S SPRMNT 22"
$ws.Range("B49").Value = "FSHRMAN'S HNY&LMN 22"
$ws.Range("B50").Value = "FSHERMAN'S BLKCRN 22"

# --- Remove the discontinued row (10036982 / KIT KAT CHO 2F 17G) ---
$ws.Rows.Item(19).Delete()
